$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 200
$ws.Range("B3").Value = 300

# Update the selected/active cell to B3
$ws.Range("B3").Select()
